$wb = $excel.ActiveWorkbook

# --- Duplicate "Sheet1" into a new sheet "Sheet2" placed right after it ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Copy the Mes / Total Despesas block (values + number formats) from Sheet1.
$ws1.Range("A1:B13").Copy($ws2.Range("A1"))

# On the new sheet the "Total Despesas" column is entered in ascending order
# (10, 20, 30, ... 120) instead of the descending order found on Sheet1.
for ($i = 0; $i -lt 12; $i++) {
    $ws2.Cells.Item($i + 2, 2).Value = ($i + 1) * 10
}

# Column B holds the same kind of content as on Sheet1 ("Total Despesas"),
# so it gets the same best-fit width.
$ws2.Columns.Item(2).ColumnWidth = $ws1.Columns.Item(2).ColumnWidth

# Turn the range into a proper table, matching the style used on the sibling
# "Despesas" sheet.
$lo = $ws2.ListObjects.Add(1, $ws2.Range("A1:B13"), $null, 1)
$lo.Name = "Table_356"
$lo.TableStyle = "Fluxo anual_2025_Despesas(6)-style"
$lo.ShowTableStyleFirstColumn = $true
$lo.ShowTableStyleLastColumn = $true
$lo.ShowAutoFilter = $false
$lo.ShowTotals = $false

# --- Update selections / active tab ---
# Sheet1 keeps the whole table selected but is no longer the active tab.
$ws1.Range("A1:B13").Select()

# Sheet2 becomes the active sheet, with the same "next empty row" selection
# Sheet1 used to have.
$ws2.Activate()
$ws2.Range("B14").Select()
